# Update COVID-19 "paises" stats workbook:
#  - refresh the "datos actualizados" timestamp
#  - refresh per-country counters for the rows whose data changed
#  - some rows changed which country occupies them because the sheet is
#    kept sorted descending by "Casos totales" (column B), so a handful
#    of countries swap places with their neighbours once the new counts
#    are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 19:08"

# Columns: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
Set-RowData 4   'Estados Unidos'         1197907 9785 179783 949113 16039 413 69011
Set-RowData 9   'Alemania'               165745  81   132700 26179  1949  0   6866
Set-RowData 20  'Ecuador'                31881   2343 3433   26879  159   5   1569
Set-RowData 21  'Suiza'                  29981   76   24500  3702   141   17  1779
Set-RowData 22  'Arabia Saudita'         28656   1645 4476   23989  143   7   191
Set-RowData 26  'Irlanda'                21772   266  13386  7067   364   16  1319
Set-RowData 37  'Rumania'                13512   349  5269   7425   243   28  818
Set-RowData 51  'Egipto'                 6813    348  1632   4745   0     7   436
Set-RowData 52  'Sudafrica'              6783    0    2549   4103   36    0   131
Set-RowData 59  'Moldavia'               4248    127  1423   2693   237   7   132
Set-RowData 77  'Islandia'               1799    0    1723   66     0     0   10
Set-RowData 99  'Principado de Andorra'  750     2    499    206    16    0   45
Set-RowData 100 'Niger'                  750     0    518    196    0     0   36
Set-RowData 137 'Martinica'              181     2    83     84     5     0   14
Set-RowData 142 'Guadalupe'              152     0    98     42     5     0   12
Set-RowData 144 'Gibraltar'              144     0    133    11     0     0   0
Set-RowData 148 'Guayana Francesa'       133     5    100    32     2     0   1
Set-RowData 152 'Suazilandia'            116     4    12     103    0     0   1
Set-RowData 153 'Trinidad yTobago'       116     0    96     12     0     0   8
Set-RowData 154 'Bermudas'               115     0    51     57     4     0   7
Set-RowData 198 'San Cristobal y Nieves' 15      0    8      7      0     0   0
Set-RowData 199 'Burundi'                15      0    7      7      0     0   1
